$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new team members to row 11: column H (FANTASTIC 9) and column J (No Mo Junk in da Trunk)
# Order matters for shared string table indices: Souvik Nath must be added before Akhilesh Paliwal
$ws.Range("J11").Value = "Souvik Nath"
$ws.Range("H11").Value = "Akhilesh Paliwal"

# Match the formatting used by the rest of the data rows (same style as H10/J10)
$ws.Range("H10").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("J10").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author left it after editing
$ws.Range("M16").Select() | Out-Null
